# Updated ligand/receptor TPM-derived values (Col4a1-Itgb1) and their
# downstream specificity / edge-weight recalculations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 93.39526366666666
$ws.Range("H2").Value = 280.185791
$ws.Range("I2").Value = 0.2167755775732346
$ws.Range("J2").Value = 0.2167755775732346
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 7199.796357654319
$ws.Range("R2").Value = 64798.16721888887
$ws.Range("S2").Value = 0.05210887478273039
$ws.Range("T2").Value = 0.05210887478273039

# Row 3
$ws.Range("G3").Value = 93.39526366666666
$ws.Range("H3").Value = 280.185791
$ws.Range("I3").Value = 0.2167755775732346
$ws.Range("J3").Value = 0.2167755775732346
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 9487.09437001651
$ws.Range("R3").Value = 85383.84933014859
$ws.Range("S3").Value = 0.06866330490772367
$ws.Range("T3").Value = 0.06866330490772365

# Row 4
$ws.Range("G4").Value = 93.39526366666666
$ws.Range("H4").Value = 280.185791
$ws.Range("I4").Value = 0.2167755775732346
$ws.Range("J4").Value = 0.2167755775732346
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 13264.62943751679
$ws.Range("R4").Value = 119381.6649376511
$ws.Range("S4").Value = 0.09600339788278055
$ws.Range("T4").Value = 0.09600339788278053

# Row 5
$ws.Range("I5").Value = 0.5566060939249745
$ws.Range("J5").Value = 0.5566060939249745
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 18486.63291571841
$ws.Range("R5").Value = 166379.6962414657
$ws.Range("S5").Value = 0.1337979009274813
$ws.Range("T5").Value = 0.1337979009274813

# Row 6
$ws.Range("I6").Value = 0.5566060939249745
$ws.Range("J6").Value = 0.5566060939249745
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1763040577195835
$ws.Range("T6").Value = 0.1763040577195835

# Row 7
$ws.Range("I7").Value = 0.5566060939249745
$ws.Range("J7").Value = 0.5566060939249745
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2465041352779098
$ws.Range("T7").Value = 0.2465041352779098

# Row 8
$ws.Range("I8").Value = 0.226618328501791
$ws.Range("J8").Value = 0.2266183285017909
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 7526.704965524492
$ws.Range("R8").Value = 67740.34468972043
$ws.Range("S8").Value = 0.05447489166247079
$ws.Range("T8").Value = 0.05447489166247078

# Row 9
$ws.Range("I9").Value = 0.226618328501791
$ws.Range("J9").Value = 0.2266183285017909
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.07178097995075257
$ws.Range("T9").Value = 0.07178097995075254

# Row 10
$ws.Range("I10").Value = 0.226618328501791
$ws.Range("J10").Value = 0.2266183285017909
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.1003624568885676
$ws.Range("T10").Value = 0.1003624568885676
